$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Julian Champagnie", "SF,PF", "San Antonio Spurs"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
